$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update values
$ws.Range("B3").Value = 0.9974789578434095
$ws.Range("C3").Value = 0.9976746246971326
$ws.Range("D3").Value = 0.9730274753669037

# Row 4 - rename GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9958856148886301
$ws.Range("C4").Value = 0.9964575028674818
$ws.Range("D4").Value = 0.9314450610255703

# Row 5 - rename AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9985529337222276
$ws.Range("C5").Value = 0.9983192370639756
$ws.Range("D5").Value = 0.9959504334631312
